$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The workers' database was re-ordered upstream, so the two account
# statement rows for "20432527 / RICARDO ANDRES YARA PACHECO" (row 16)
# and "1010197164 / DIANA MARIA PAYARES PEREZ" (row 18) swap places
# (document id, name, period and overdue-value columns), while the
# middle row (45532902 / SISSY EMPERATRIZ ALGARIN MENDOZA) is untouched.

# Row 16: 20432527 / RICARDO ANDRES YARA PACHECO / 2305 / 10000
#   ->    1010197164 / DIANA MARIA PAYARES PEREZ / 2404 / 8000
$ws.Range("C16").Value = "1010197164"
$ws.Range("D16").Value = "DIANA MARIA PAYARES PEREZ"
$ws.Range("E16").Value = "2404"
$ws.Range("F16").Value = 8000

# Row 18: 1010197164 / DIANA MARIA PAYARES PEREZ / 2404 / 8000
#   ->    20432527 / RICARDO ANDRES YARA PACHECO / 2305 / 10000
$ws.Range("C18").Value = "20432527"
$ws.Range("D18").Value = "RICARDO ANDRES YARA PACHECO"
$ws.Range("E18").Value = "2305"
$ws.Range("F18").Value = 10000

# Re-fit the data columns to the (now different) content widths, as
# Excel does automatically for bestFit columns when their text changes.
$ws.Range("B15:J18").Columns.AutoFit() | Out-Null
